$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The pinout list had a duplicated "RAM_A5" entry (rows 193 & 194 both read
# "RAM_A5"); remove the stray duplicate row so the list shifts back in line.
$ws.Rows(194).Delete()

# Flag the four GND pins that sit between connector blocks with a "-----"
# separator note in column B (typed with a leading apostrophe so Excel
# stores it as literal text rather than trying to parse it).
$ws.Range("B52").Value = "'-----"
$ws.Range("B104").Value = "'-----"
$ws.Range("B156").Value = "'-----"
$ws.Range("B208").Value = "'-----"

# Fix the typo'd pin name (missing underscore).
$ws.Range("A183").Value = "RAM_UCASU#"
